$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'73.378.39"
$ws.Range("E2").Value = "  +1.82%  "

# Row 3
$ws.Range("D3").Value = "'3.986.16"
$ws.Range("E3").Value = "  -1.13%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").Value = "'615.36"
$ws.Range("E5").Value = "  +14.21%  "

# Row 6
$ws.Range("D6").Value = "'166.26"
$ws.Range("E6").Value = "  +11.66%  "

# Row 7
$ws.Range("D7").Value = "'0.683"
$ws.Range("E7").Value = "  -1.68%  "

# Row 8
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").Value = "'0.754"
$ws.Range("E9").Value = "  +0.51%  "

# Row 10
$ws.Range("D10").Value = "'0.186"
$ws.Range("E10").Value = "  +8.04%  "

# Row 11
$ws.Range("D11").Value = "'56.43"
$ws.Range("E11").Value = "  +6.24%  "

# Row 12
$ws.Range("D12").Value = "'0.0000339"
$ws.Range("E12").Value = "  +2.98%  "

# Row 13
$ws.Range("D13").Value = "'11.09"
$ws.Range("E13").Value = "  +2.09%  "

# Row 14
$ws.Range("D14").Value = "'4.626.51"
$ws.Range("E14").Value = "  -1.19%  "

# Row 15
$ws.Range("D15").Value = "'3.998.93"
$ws.Range("E15").Value = "  -1.34%  "

# Row 16
$ws.Range("D16").Value = "'1.24"
$ws.Range("E16").Value = "  +3.37%  "

# Row 17
$ws.Range("D17").Value = "'14.15"
$ws.Range("E17").Value = "  -0.77%  "

# Row 18
$ws.Range("D18").Value = "'20.51"
$ws.Range("E18").Value = "  -0.61%  "

# Row 19
$ws.Range("D19").Value = "'73.282.68"
$ws.Range("E19").Value = "  +1.62%  "

# Row 20
$ws.Range("E20").Value = "  -0.33%  "

# Row 21
$ws.Range("D21").Value = "'440.54"
$ws.Range("E21").Value = "  +0.44%  "

# Row 22
$ws.Range("D22").Value = "'4.90"
$ws.Range("E22").Value = "  +14.73%  "

# Row 23
$ws.Range("D23").Value = "'95.78"
$ws.Range("E23").Value = "  -2.16%  "

# Row 24
$ws.Range("D24").Value = "'3.36"
$ws.Range("E24").Value = "  -3.81%  "

# Row 25
$ws.Range("D25").Value = "'14.14"
$ws.Range("E25").Value = "  -2.97%  "

# Row 26
$ws.Range("D26").Value = "'4.08"
$ws.Range("E26").Value = "  -8.92%  "

# Row 27
$ws.Range("D27").Value = "'11.12"
$ws.Range("E27").Value = "  -1.12%  "

# Row 28
$ws.Range("E28").Value = "  -0.42%  "

# Row 29
$ws.Range("D29").Value = "'10.46"
$ws.Range("E29").Value = "  -1.85%  "

# Row 30
$ws.Range("D30").Value = "'36.12"
$ws.Range("E30").Value = "  -2.60%  "

# Row 31
$ws.Range("D31").Value = "'7.81"
$ws.Range("E31").Value = "  -6.39%  "

# Row 32
$ws.Range("D32").Value = "'13.68"
$ws.Range("E32").Value = "  +1.11%  "

# Row 33
$ws.Range("E33").Value = "  -2.79%  "

# Row 34
$ws.Range("D34").Value = "'71.28"
$ws.Range("E34").Value = "  +6.74%  "

# Row 35
$ws.Range("D35").Value = "'47.58"
$ws.Range("E35").Value = "  -3.42%  "

# Row 36
$ws.Range("D36").Value = "'0.0000100"
$ws.Range("E36").Value = "  +10.27%  "

# Row 37
$ws.Range("D37").Value = "'637.77"
$ws.Range("E37").Value = "  -6.33%  "

# Row 38
$ws.Range("D38").Value = "'0.431"
$ws.Range("E38").Value = "  -5.72%  "

# Row 39
$ws.Range("D39").Value = "'3.42"
$ws.Range("E39").Value = "  +0.73%  "

# Row 40
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.08%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.146"
$ws.Range("E41").Value = "  -1.25%  "

# Row 42
$ws.Range("D42").Value = "'11.03"
$ws.Range("E42").Value = "  -2.75%  "

# Row 43
$ws.Range("E43").Value = "  +0.19%  "

# Row 44
$ws.Range("D44").Value = "'3.28"
$ws.Range("E44").Value = "  -3.78%  "

# Row 45
$ws.Range("E45").Value = "  -1.55%  "

# Row 46
$ws.Range("E46").Value = "  -0.78%  "

# Row 47
$ws.Range("D47").Value = "'3.42"
$ws.Range("E47").Value = "  +2.43%  "

# Row 48
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.90"
$ws.Range("E48").Value = "  +27.62%  "

# Row 49
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").Value = "'2.61"
$ws.Range("E49").Value = "  -0.94%  "

# Row 50
$ws.Range("D50").Value = "'2.858.69"
$ws.Range("E50").Value = "  +3.10%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'3.02"
$ws.Range("E51").Value = "  -3.25%  "
